$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 3 already has B3="Anfangs Code", C3=56.2
# Add the new name column (A) and new rows (4-6)

$ws.Range("A3").Value = "Janes"
$ws.Range("B3").Value = "Anfangs Code"
$ws.Range("C3").Value = 56.2

$ws.Range("A4").Value = "Elias"
$ws.Range("B4").Value = "Variablen Char"
$ws.Range("C4").Value = 30

$ws.Range("A5").Value = "Elias"
$ws.Range("B5").Value = "array => spielfeld"
$ws.Range("C5").Value = 30

$ws.Range("A6").Value = "Janes"
$ws.Range("B6").Value = "zaehleLebende vereinfacht"
$ws.Range("C6").Value = 27

$ws.Range("D9").Select()
